$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto price/volume refresh (GitHub Actions bot).
# Column D = Price (text), Column E = Volume(1h) change (text, padded).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.801.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.526.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.68"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.51"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.109"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.923.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.479.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.862.56"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.61"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.29"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.05"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.37"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "41.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.43"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.66"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.36"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0789"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.10%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -12.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0305"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.80"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.003.94"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.13"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.09"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.89"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.776.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.25%  "
